# Generate Report for Handback
# Updates handback-status timestamps and status text to reflect a newer
# report generation run.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: Latest HO Xliff Generate Date ---
$overview.Range("G2").Value = "2016-08-19 00:16:57"
$overview.Range("G3").Value = "2016-08-19 00:16:57"

# --- zh-cn sheet ---
# Priority changed from "ht" (human translation) to "mt" (machine translation)
$zhcn.Range("E2").Value = "mt"
$zhcn.Range("E3").Value = "mt"

# Correspond Handoff Datetime
$zhcn.Range("H2").Value = "2016-08-19 00:16:51"
$zhcn.Range("H3").Value = "2016-08-19 00:16:51"

# Correspond Handback DateTime
$zhcn.Range("K2").Value = "2016-08-19 00:17:15"
$zhcn.Range("K3").Value = "2016-08-19 00:17:15"

# --- de-de sheet ---
# Priority changed from "ht" to "mt"
$dede.Range("E2").Value = "mt"
$dede.Range("E3").Value = "mt"

# Correspond Handback DateTime
$dede.Range("K2").Value = "2016-08-19 00:17:22"
$dede.Range("K3").Value = "2016-08-19 00:17:22"
